$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row labels (B column, rows 3-29), replacing the old set of labels (rows 3-19)
$labels = @(
    "Spiral5",
    "RotRing OmegaMax-90",
    "Equal Angle",
    "Tilt Rotate",
    "CLR",
    "Rizzie Hex",
    "Thomas Hex",
    "Tilt Rotate_Partial",
    "RotRing OmegaMax-60",
    "Equal Angle_Partial",
    "Rizzie Hex_Partial",
    "ND Single",
    "RD Single",
    "TD Single",
    "Morris Single",
    "Ring Perpendicular to ND",
    "Ring Perpendicular to RD",
    "Ring Perpendicular to TD",
    "OffsetFTD",
    "OffsetATD",
    "OffsetF45",
    "OffsetA45",
    "OffsetFRD",
    "OffsetARD",
    "Gaussian Quadrature",
    "Michael-CCHex",
    "Michael-SNHex"
)

# New column headers (C2:W2), replacing the previous set
$headers = @(
    "[4, 2, 2]",
    "[5, 1, 1]",
    "[2, 2, 2]",
    "[1, 1, 1]",
    "[3, 1, 1]",
    "[3, 3, 1]",
    "[2, 2, 0]",
    "[2, 0, 0]",
    "[3, 3, 3]",
    "[4, 0, 0]",
    "[4, 2, 0]",
    "1Pair-A",
    "1Pair-B",
    "2Pairs-A",
    "2Pairs-B",
    "3Pairs-A",
    "3Pairs-B",
    "3Pairs-C",
    "4Pairs",
    "5A4F",
    "MaxUnique"
)

# Clear the old grid entirely (old extent was A1:AQ19) and rebuild fresh
$ws.Range("A1:AQ19").Clear()

# Row 1: sequence numbers 0..21 in B1:W1
for ($i = 0; $i -le 21; $i++) {
    $cell = $ws.Cells.Item(1, 2 + $i)
    $cell.Value = $i
    $cell.Style = "Normal"
}
$ws.Range("B1:W1").Font.Bold = $true

# Row 2: A2=0, B2="HKL", C2:W2 = headers
$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(2, 2).Value = "HKL"
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(2, 3 + $i).Value = $headers[$i]
}

# Rows 3-29: A = row index (row-2), B = label, C:W = 1
for ($r = 0; $r -lt $labels.Length; $r++) {
    $row = 3 + $r
    $ws.Cells.Item($row, 1).Value = $r + 1
    $ws.Cells.Item($row, 2).Value = $labels[$r]
    for ($c = 3; $c -le 23; $c++) {
        $ws.Cells.Item($row, $c).Value = 1
    }
}

$ws.Range("A2:A29").Font.Bold = $true
